# Update "want to go" counts (column F) on the "展览" sheet and the
# aggregated "全部类型" sheet to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F value }
$sheetUpdates = @{
    "展览" = @{
        2  = 3125
        3  = 523
        4  = 1078
        6  = 27
        8  = 30
        9  = 1108
        10 = 15492
        11 = 227
        14 = 6111
        16 = 103
        19 = 105
        21 = 27
        22 = 114
        23 = 8
        27 = 855
        28 = 18
        30 = 138
        31 = 10981
        37 = 261
    }
    "全部类型" = @{
        3  = 3125
        4  = 523
        5  = 1078
        7  = 27
        9  = 30
        10 = 1108
        11 = 15492
        12 = 227
        15 = 6111
        17 = 103
        20 = 105
        22 = 27
        23 = 114
        24 = 8
        28 = 855
        29 = 18
        31 = 138
        33 = 10981
        39 = 261
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
